$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.659.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.564.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.97%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.86'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.91%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0586'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.788.14'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.569.03'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.667.86'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.517'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.09%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.47'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.58'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.33'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0680'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.92'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.04'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.90%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.23'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0458'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.01%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.99%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.401.58'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.98'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.02'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.46'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.71'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.77%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.88%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.518'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.94'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.768'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0458'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.87'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.22'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.700.19'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.868'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.72'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.28'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.86%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.51%  '
